$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = -0.7091679257846845
    "C2"  = -0.2116494899511859
    "D2"  = -0.8069708324279797
    "E2"  = -0.442835019116573
    "F2"  = -0.02123990589338218
    "G2"  = -0.1707229344526847
    "H2"  = -0.1169158257196027

    "B3"  = -0.5825132856353434
    "C3"  = -1.177834628112137
    "D3"  = -0.8136988148007305
    "E3"  = -0.3921037015775397
    "F3"  = -0.5415867301368422
    "G3"  = -0.4877796214037602

    "B4"  = -0.6994291822578886
    "C4"  = -0.3352933689464819
    "D4"  = 0.08630174427670884
    "E4"  = -0.06318128428259367
    "F4"  = -0.009374175549511699
    "G4"  = -0.334123210679779
    "H4"  = 0.001243794298757539
    "I4"  = 0.09478491567901273
    "J4"  = -0.3880619836346724

    "B5"  = -0.1970460893006987
    "C5"  = 0.2245490239224921
    "D5"  = 0.07506599536318959
    "E5"  = 0.1288731040962716
    "F5"  = -0.1958759310339957
    "G5"  = 0.1394910739445408
    "H5"  = 0.233032195324796
    "I5"  = -0.2498147039888892

    "B6"  = 0.1735200387801499
    "C6"  = 0.02403701022084742
    "D6"  = 0.0778441189539294
    "E6"  = -0.2469049161763379
    "F6"  = 0.08846208880219864
    "G6"  = 0.1820032101824538
    "H6"  = -0.3008436891312313

    "B7"  = -0.009228378467612708
    "C7"  = 0.04457873026546927
    "D7"  = -0.280170304864798
    "E7"  = 0.05519670011373851
    "F7"  = 0.1487378214939937
    "G7"  = -0.3341090778196915

    "B8"  = 0.1858859418247861
    "C8"  = -0.1388630933054812
    "D8"  = 0.1965039116730553
    "E8"  = 0.2900450330533105
    "F8"  = -0.1928018662603747
    "G8"  = -0.08321328366289282
    "H8"  = -0.05405384115675094
    "I8"  = -0.2527483093463515

    "B9"  = -0.2674987412282675
    "C9"  = 0.067868263750269
    "D9"  = 0.1614093851305242
    "E9"  = -0.321437514183161
    "F9"  = -0.2118489315856791
    "G9"  = -0.1826894890795372
    "H9"  = -0.3813839572691378

    "B10" = 0.07761257013415129
    "C10" = 0.1711536915144065
    "D10" = -0.3116932077992787
    "E10" = -0.2021046252017968
    "F10" = -0.172945182695655
    "G10" = -0.3716396508852555

    "B11" = -0.1119472883713244
    "C11" = -0.5947941876850096
    "D11" = -0.4852056050875277
    "E11" = -0.4560461625813858
    "F11" = -0.6547406307709864

    "B12" = -0.2799851348954449
    "C12" = -0.170396552297963
    "D12" = -0.1412371097918211
    "E12" = -0.3399315779814217

    "B13" = -0.1408857896223169
    "C13" = -0.111726347116175
    "D13" = -0.3104208153057756

    "B14" = -0.2487487154747626
    "C14" = -0.4474431836643632

    "B15" = -0.3133716091653966
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
